$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values (e.g. "1.00", "0.997",
# "65.317.19") as plain text in the source workbook. Force Text format on just
# the D cells we are about to rewrite so Excel does not auto-coerce them into
# real numbers (which would silently normalize "1.00" -> 1, drop the thousands-
# dot grouping in values like "65.317.19", etc.). The "Volume(1h)" column (E)
# always contains a "%" sign plus padding spaces, so it is never parsed as a
# number and needs no special handling.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.317.19'
$ws.Range("E2").Value = '  +2.46%  '
$ws.Range("D3").Value = '3.443.12'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '578.77'
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("D6").Value = '167.25'
$ws.Range("E6").Value = '  +6.38%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.443.14'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("D9").Value = '0.561'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = '7.21'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").Value = '0.122'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").Value = '0.428'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '4.046.05'
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '27.41'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("D17").Value = '65.323.63'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = '3.450.59'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '6.20'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '13.75'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '380.83'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '7.90'
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = '71.47'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").Value = '0.519'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +1.73%  '
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  +3.48%  '
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '6.30'
$ws.Range("E30").Value = '  +5.38%  '
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  +3.04%  '
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("D33").Value = '23.12'
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '7.27'
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("D36").Value = '1.51'
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("D37").Value = '160.26'
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '0.895'
$ws.Range("E38").Value = '  +9.74%  '
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").Value = '0.0737'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.806.86'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '26.01'
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").Value = '6.58'
$ws.Range("E43").Value = '  +3.06%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '42.98'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '26.48'
$ws.Range("E45").Value = '  +4.68%  '
$ws.Range("D46").Value = '4.43'
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0307'
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.49'
$ws.Range("E48").Value = '  +6.46%  '
$ws.Range("D49").Value = '347.22'
$ws.Range("E49").Value = '  +6.75%  '
$ws.Range("D50").Value = '1.06'
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").Value = '32.12'
$ws.Range("E51").Value = '  +7.42%  '
